# Update countries & provincias Spain
# - refresh the "last updated" timestamp
# - refresh case counts for several countries whose daily figures changed
# - a few countries leap-frogged their neighbour in the ranking, so the
#   row that used to hold the lower-ranked country now gets the
#   higher-ranked country's (freshly updated) figures, and the
#   previously top row's *old* figures slide down to the neighbour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 14:22"

# --- simple value refreshes (no reordering) ---------------------------
# Portugal (row 18)
$ws.Range("B18").Value = 15472
$ws.Range("C18").Value = 1516
$ws.Range("D18").Value = 233
$ws.Range("E18").Value = 14804
$ws.Range("F18").Value = 226
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = 435

# Austria (row 19)
$ws.Range("B19").Value = 13453
$ws.Range("C19").Value = 209
$ws.Range("E19").Value = 7070

# Suecia (row 23)
$ws.Range("B23").Value = 9685
$ws.Range("C23").Value = 544
$ws.Range("E23").Value = 8610
$ws.Range("F23").Value = 749
$ws.Range("G23").Value = 77
$ws.Range("H23").Value = 870

# --- Polonia / Dinamarca swap places (rows 29-30) ----------------------
# Dinamarca overtakes Polonia and gets fresh numbers; Polonia keeps its
# previous (unchanged) figures but now sits one row lower.
$ws.Range("A29").Value = "Dinamarca"
$ws.Range("B29").Value = 5819
$ws.Range("C29").Value = 184
$ws.Range("D29").Value = 1773
$ws.Range("E29").Value = 3799
$ws.Range("F29").Value = 113
$ws.Range("G29").Value = 10
$ws.Range("H29").Value = 247

$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 5742
$ws.Range("C30").Value = 167
$ws.Range("D30").Value = 318
$ws.Range("E30").Value = 5249
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 175

# --- Marruecos / Croacia swap places (rows 60-61) ----------------------
$ws.Range("A60").Value = "Croacia"
$ws.Range("B60").Value = 1495
$ws.Range("C60").Value = 88
$ws.Range("D60").Value = 231
$ws.Range("E60").Value = 1243
$ws.Range("F60").Value = 34
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 21

$ws.Range("A61").Value = "Marruecos"
$ws.Range("B61").Value = 1431
$ws.Range("C61").Value = 57
$ws.Range("D61").Value = 114
$ws.Range("E61").Value = 1212
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = 8
$ws.Range("H61").Value = 105

# --- Eslovaquia / Republica de Macedonia swap places (rows 78-79) ------
$ws.Range("A78").Value = "Republica de Macedonia"
$ws.Range("B78").Value = 711
$ws.Range("C78").Value = 48
$ws.Range("D78").Value = 41
$ws.Range("E78").Value = 638
$ws.Range("F78").Value = 15
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 32

$ws.Range("A79").Value = "Eslovaquia"
$ws.Range("B79").Value = 701
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 23
$ws.Range("E79").Value = 676
$ws.Range("F79").Value = 5
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 2

# --- Sri Lanka / Isla de Man / Mayotte reorder (rows 115-117) ----------
# Mayotte jumps to the top of this trio with fresh numbers; Sri Lanka and
# Isla de Man both keep their previous figures but slide one row down.
$ws.Range("A115").Value = "Mayotte"
$ws.Range("B115").Value = 191
$ws.Range("C115").Value = 7
$ws.Range("D115").Value = 50
$ws.Range("E115").Value = 139
$ws.Range("F115").Value = 4
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 2

$ws.Range("A116").Value = "Sri Lanka"
$ws.Range("B116").Value = 190
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 50
$ws.Range("E116").Value = 133
$ws.Range("F116").Value = 5
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 7

$ws.Range("A117").Value = "Isla de Man"
$ws.Range("B117").Value = 190
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 92
$ws.Range("E117").Value = 97
$ws.Range("F117").Value = 12
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 1
